$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "FAILED"
$ws.Range("E2").Value = 8.488922890999675
$ws.Range("F2").Value = "2022-08-22T17:37:53"

$message = @"
selenium.common.exceptions.InvalidArgumentException: Message: invalid argument: 'text' is empty
  (Session info: chrome=104.0.5112.101)
Stacktrace:
0   chromedriver                        0x000000010ece4149 chromedriver + 4469065
1   chromedriver                        0x000000010ec6e233 chromedriver + 3985971
2   chromedriver                        0x000000010e904fe8 chromedriver + 409576
3   chromedriver                        0x000000010e933a78 chromedriver + 600696
4   chromedriver                        0x000000010e958c62 chromedriver + 752738
5   chromedriver                        0x000000010e92fe35 chromedriver + 585269
6   chromedriver                        0x000000010e958d6e chromedriver + 753006
7   chromedriver                        0x000000010e96b611 chromedriver + 828945
8   chromedriver                        0x000000010e958b53 chromedriver + 752467
9   chromedriver                        0x000000010e92e905 chromedriver + 579845
10  chromedriver                        0x000000010e92f955 chromedriver + 584021
11  chromedriver                        0x000000010ecb56ad chromedriver + 4277933
12  chromedriver                        0x000000010ecb9b3a chromedriver + 4295482
13  chromedriver                        0x000000010ecbecdf chromedriver + 4316383
14  chromedriver                        0x000000010ecba857 chromedriver + 4298839
15  chromedriver                        0x000000010ec9364f chromedriver + 4138575
16  chromedriver                        0x000000010ecd51f8 chromedriver + 4407800
17  chromedriver                        0x000000010ecd537f chromedriver + 4408191
18  chromedriver                        0x000000010ecebcb5 chromedriver + 4500661
19  libsystem_pthread.dylib             0x00007ff80a4144e1 _pthread_start + 125
20  libsystem_pthread.dylib             0x00007ff80a40ff6b thread_start + 15
"@

$ws.Range("G2").Value = $message
